$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 03:11"

# --- Refresh totals for existing countries (no reorder) ---
# Estados Unidos (row 4)
$ws.Range("B4").Value = 6290425
$ws.Range("C4").Value = 40899
$ws.Range("D4").Value = 3546627
$ws.Range("E4").Value = 2553857
$ws.Range("G4").Value = 1067
$ws.Range("H4").Value = 189941

# Venezuela (row 56)
$ws.Range("B56").Value = 48883
$ws.Range("C56").Value = 1127
$ws.Range("D56").Value = 39912
$ws.Range("E56").Value = 8573
$ws.Range("G56").Value = 7
$ws.Range("H56").Value = 398

# Guinea Ecuatorial (row 111)
$ws.Range("D111").Value = 4390
$ws.Range("E111").Value = 492

# --- Insert "Congo" as a new row right after Suazilandia (row 115), pushing
#     Ruanda..Eslovaquia down by one row and dropping Congo's old row (123) ---
$ws.Range("A116").Value = "Congo"
$ws.Range("B116").Value = 4628
$ws.Range("C116").Value = 649
$ws.Range("D116").Value = 1742
$ws.Range("E116").Value = 2784
$ws.Range("G116").Value = 24
$ws.Range("H116").Value = 102

$ws.Range("A117").Value = "Ruanda"
$ws.Range("B117").Value = 4218
$ws.Range("C117").Value = 76
$ws.Range("D117").Value = 2071
$ws.Range("E117").Value = 2130
$ws.Range("H117").Value = 17

$ws.Range("A118").Value = "Tunez"
$ws.Range("B118").Value = 4196
$ws.Range("C118").Value = 233
$ws.Range("D118").Value = 1628
$ws.Range("E118").Value = 2487
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 81

$ws.Range("A119").Value = "Surinam"
$ws.Range("B119").Value = 4149
$ws.Range("C119").Value = 60
$ws.Range("D119").Value = 3272
$ws.Range("E119").Value = 805
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 72

$ws.Range("A120").Value = "Cuba"
$ws.Range("B120").Value = 4126
$ws.Range("C120").Value = 61
$ws.Range("D120").Value = 3458
$ws.Range("E120").Value = 570
$ws.Range("G120").Value = 3
$ws.Range("H120").Value = 98

$ws.Range("A121").Value = "Mozambique"
$ws.Range("B121").Value = 4117
$ws.Range("D121").Value = 2170
$ws.Range("E121").Value = 1922
$ws.Range("G121").Value = 2
$ws.Range("H121").Value = 25

$ws.Range("A122").Value = "Cabo Verde"
$ws.Range("B122").Value = 4048
$ws.Range("C122").Value = 78
$ws.Range("D122").Value = 3460
$ws.Range("E122").Value = 547
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 41

$ws.Range("A123").Value = "Eslovaquia"
$ws.Range("B123").Value = 4042
$ws.Range("C123").Value = 53
$ws.Range("D123").Value = 2523
$ws.Range("E123").Value = 1486
$ws.Range("H123").Value = 33

# --- Bahamas refresh (row 138) ---
$ws.Range("B138").Value = 2337
$ws.Range("C138").Value = 61
$ws.Range("D138").Value = 849
$ws.Range("E138").Value = 1438

# --- Move "Guadalupe" up ahead of "Liberia" (rows 158-159) ---
$ws.Range("A158").Value = "Guadalupe"
$ws.Range("B158").Value = 1363
$ws.Range("C158").Value = 94
$ws.Range("D158").Value = 336
$ws.Range("E158").Value = 1009
$ws.Range("G158").Value = 2
$ws.Range("H158").Value = 18

$ws.Range("A159").Value = "Liberia"
$ws.Range("B159").Value = 1305
$ws.Range("D159").Value = 1162
$ws.Range("E159").Value = 61
$ws.Range("H159").Value = 82

# Martinica (row 168)
$ws.Range("B168").Value = 754
$ws.Range("C168").Value = 7
$ws.Range("E168").Value = 638

# San Marino (row 169)
$ws.Range("B169").Value = 735
$ws.Range("C169").Value = 20
$ws.Range("E169").Value = 33

# Monaco (row 192)
$ws.Range("B192").Value = 142
$ws.Range("C192").Value = 2
$ws.Range("D192").Value = 91
$ws.Range("E192").Value = 50
